$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q1" sheet right before the "总计" (total) sheet ---
# NOTE: worksheet object handles in this runtime are resolved by *position*,
# not a stable identity. Once Worksheets.Add() shifts positions around, a
# handle obtained beforehand (e.g. $wb.Worksheets.Item("总计")) can silently
# start referring to a different sheet. So: grab "总计" only to anchor the
# Add() call, then re-resolve it by name again afterwards before touching it.
$totalSheetBeforeAdd = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBeforeAdd)
$newSheet.Name = "2022-Q1"

# Copy header (row1) + first-data-row (row2) formatting from the "2021-Q4" sheet,
# which already has the fund-holdings layout we need to replicate.
$srcFmt = $wb.Worksheets.Item("2021-Q4")
$srcFmt.Range("A1:H2").Copy() | Out-Null
$newSheet.Range("A1:H2").PasteSpecial(-4122) | Out-Null

# Stamp that same data-row formatting down across all 13 data rows (rows 2-14).
$newSheet.Range("A2:H2").Copy() | Out-Null
$newSheet.Range("A3:H14").PasteSpecial(-4122) | Out-Null

# Header row text
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Fund-holdings data rows (A = running index, H = numeric rank; the rest are text)
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'002121"
$newSheet.Cells.Item(2,3).Value = "广发沪港深新起点股票A"
$newSheet.Cells.Item(2,4).Value = "'34.56"
$newSheet.Cells.Item(2,5).Value = "'91.46"
$newSheet.Cells.Item(2,6).Value = "'6.16"
$newSheet.Cells.Item(2,7).Value = "'2.1289"
$newSheet.Cells.Item(2,8).Value = 5
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'009265"
$newSheet.Cells.Item(3,3).Value = "易方达消费精选股票"
$newSheet.Cells.Item(3,4).Value = "'50.44"
$newSheet.Cells.Item(3,5).Value = "'82.26"
$newSheet.Cells.Item(3,6).Value = "'3.58"
$newSheet.Cells.Item(3,7).Value = "'1.8058"
$newSheet.Cells.Item(3,8).Value = 10
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'009896"
$newSheet.Cells.Item(4,3).Value = "广发港股通成长精选股票A"
$newSheet.Cells.Item(4,4).Value = "'27.73"
$newSheet.Cells.Item(4,5).Value = "'89.63"
$newSheet.Cells.Item(4,6).Value = "'6.31"
$newSheet.Cells.Item(4,7).Value = "'1.7498"
$newSheet.Cells.Item(4,8).Value = 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'070012"
$newSheet.Cells.Item(5,3).Value = "嘉实海外中国混合(QDII)"
$newSheet.Cells.Item(5,4).Value = "'18.97"
$newSheet.Cells.Item(5,5).Value = "'89.48"
$newSheet.Cells.Item(5,6).Value = "'3.36"
$newSheet.Cells.Item(5,7).Value = "'0.6374"
$newSheet.Cells.Item(5,8).Value = 10
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'010350"
$newSheet.Cells.Item(6,3).Value = "景顺长城品质长青混合"
$newSheet.Cells.Item(6,4).Value = "'13.35"
$newSheet.Cells.Item(6,5).Value = "'92.33"
$newSheet.Cells.Item(6,6).Value = "'4.19"
$newSheet.Cells.Item(6,7).Value = "'0.5594"
$newSheet.Cells.Item(6,8).Value = 10
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'009897"
$newSheet.Cells.Item(7,3).Value = "广发港股通成长精选股票C"
$newSheet.Cells.Item(7,4).Value = "'6.49"
$newSheet.Cells.Item(7,5).Value = "'89.63"
$newSheet.Cells.Item(7,6).Value = "'6.31"
$newSheet.Cells.Item(7,7).Value = "'0.4095"
$newSheet.Cells.Item(7,8).Value = 5
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'006752"
$newSheet.Cells.Item(8,3).Value = "天弘港股通精选灵活配置混合A"
$newSheet.Cells.Item(8,4).Value = "'6.86"
$newSheet.Cells.Item(8,5).Value = "'85.37"
$newSheet.Cells.Item(8,6).Value = "'5.45"
$newSheet.Cells.Item(8,7).Value = "'0.3739"
$newSheet.Cells.Item(8,8).Value = 4
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'006753"
$newSheet.Cells.Item(9,3).Value = "天弘港股通精选灵活配置混合C"
$newSheet.Cells.Item(9,4).Value = "'2.49"
$newSheet.Cells.Item(9,5).Value = "'85.37"
$newSheet.Cells.Item(9,6).Value = "'5.45"
$newSheet.Cells.Item(9,7).Value = "'0.1357"
$newSheet.Cells.Item(9,8).Value = 4
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'013009"
$newSheet.Cells.Item(10,3).Value = "万家港股通精选混合A"
$newSheet.Cells.Item(10,4).Value = "'2.97"
$newSheet.Cells.Item(10,5).Value = "'81.62"
$newSheet.Cells.Item(10,6).Value = "'3.86"
$newSheet.Cells.Item(10,7).Value = "'0.1146"
$newSheet.Cells.Item(10,8).Value = 9
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "'010024"
$newSheet.Cells.Item(11,3).Value = "广发沪港深新起点股票C"
$newSheet.Cells.Item(11,4).Value = "'0.72"
$newSheet.Cells.Item(11,5).Value = "'91.46"
$newSheet.Cells.Item(11,6).Value = "'6.16"
$newSheet.Cells.Item(11,7).Value = "'0.0444"
$newSheet.Cells.Item(11,8).Value = 5
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "'013010"
$newSheet.Cells.Item(12,3).Value = "万家港股通精选混合C"
$newSheet.Cells.Item(12,4).Value = "'0.86"
$newSheet.Cells.Item(12,5).Value = "'81.62"
$newSheet.Cells.Item(12,6).Value = "'3.86"
$newSheet.Cells.Item(12,7).Value = "'0.0332"
$newSheet.Cells.Item(12,8).Value = 9
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "'009733"
$newSheet.Cells.Item(13,3).Value = "创金合信港股通大消费精选股票A"
$newSheet.Cells.Item(13,4).Value = "'0.13"
$newSheet.Cells.Item(13,5).Value = "'82.28"
$newSheet.Cells.Item(13,6).Value = "'6.82"
$newSheet.Cells.Item(13,7).Value = "'0.0089"
$newSheet.Cells.Item(13,8).Value = 2
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "'009734"
$newSheet.Cells.Item(14,3).Value = "创金合信港股通大消费精选股票C"
$newSheet.Cells.Item(14,4).Value = "'0.07"
$newSheet.Cells.Item(14,5).Value = "'82.28"
$newSheet.Cells.Item(14,6).Value = "'6.82"
$newSheet.Cells.Item(14,7).Value = "'0.0048"
$newSheet.Cells.Item(14,8).Value = 2


# --- 2. Insert the new "2022-Q1" summary row at the top of the "总计" sheet's data ---
# Re-resolve "总计" by name now that the sheet collection has changed shape.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Re-apply the bold/bordered "index" style (column A) that the other data rows use,
# copying it from the row directly below (the former row 2, now shifted to row 3).
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 13
$totalSheet.Cells.Item(2,4).Value = 8.01

# Renumber the running index in column A for the rows that got pushed down
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
